$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $fullText = $tr.Text
        if ($fullText -ne $null -and $fullText.Contains("+4+6=11")) {
            $idx = $fullText.IndexOf("+4+6=11")
            $sub = $tr.Characters($idx + 1, 7)
            $sub.Text = "+4+5=11"
        }
    }
}
